$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16088714.006621
$ws.Range("D2").Value = 44.245455

$ws.Range("B3").Value = 4388189.064076
$ws.Range("D3").Value = 6.033963
$ws.Range("E3").Value = 0.002805

$ws.Range("B4").Value = 81088175.76444501
$ws.Range("C4").Value = 223

$ws.Range("G5").Value = -173.405502
$ws.Range("H5").Value = -436.762856
$ws.Range("I5").Value = 89.951851
$ws.Range("J5").Value = 0.26812

$ws.Range("G6").Value = 138.727043
$ws.Range("H6").Value = -144.853064
$ws.Range("I6").Value = 422.307151
$ws.Range("J6").Value = 0.481862

$ws.Range("G7").Value = 312.132546
$ws.Range("H7").Value = 97.74303500000001
$ws.Range("I7").Value = 526.522056
$ws.Range("J7").Value = 0.002033
